# Append DRS review rows for Match 24 (RR vs GT) to the data table.
# Columns: A Match, B Home Team, C Away Team, D Innings, E Batting Team,
#          F Fielding Team, G Over, H Review By, I Umpire, J Umpire Abbreviation,
#          K Decision Challenged, L Original Decision, M DRS Decision,
#          N Batter, O Bowler, P Result, Q Umpires Call

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75
$ws.Cells.Item(75,1).Value = 24
$ws.Cells.Item(75,2).Value = "RR"
$ws.Cells.Item(75,3).Value = "GT"
$ws.Cells.Item(75,4).Value = 1
$ws.Cells.Item(75,5).Value = "RR"
$ws.Cells.Item(75,6).Value = "GT"
$ws.Cells.Item(75,7).Value = 17
$ws.Cells.Item(75,8).Value = "RR"
$ws.Cells.Item(75,9).Value = "Vinod Seshan"
$ws.Cells.Item(75,10).Value = "VS"
$ws.Cells.Item(75,11).Value = "Wide"
$ws.Cells.Item(75,12).Value = "Not Called"
$ws.Cells.Item(75,13).Value = "Called"
$ws.Cells.Item(75,14).Value = "SV Samson"
$ws.Cells.Item(75,15).Value = "MM Sharma"
$ws.Cells.Item(75,16).Value = "Successful"
$ws.Cells.Item(75,17).Value = "No"

# Row 76
$ws.Cells.Item(76,1).Value = 24
$ws.Cells.Item(76,2).Value = "RR"
$ws.Cells.Item(76,3).Value = "GT"
$ws.Cells.Item(76,4).Value = 1
$ws.Cells.Item(76,5).Value = "RR"
$ws.Cells.Item(76,6).Value = "GT"
$ws.Cells.Item(76,7).Value = 17
$ws.Cells.Item(76,8).Value = "GT"
$ws.Cells.Item(76,9).Value = "Vinod Seshan"
$ws.Cells.Item(76,10).Value = "VS"
$ws.Cells.Item(76,11).Value = "Wide"
$ws.Cells.Item(76,12).Value = "Called"
$ws.Cells.Item(76,13).Value = "Called"
$ws.Cells.Item(76,14).Value = "SV Samson"
$ws.Cells.Item(76,15).Value = "MM Sharma"
$ws.Cells.Item(76,16).Value = "Unsuccessful"
$ws.Cells.Item(76,17).Value = "No"

# Row 77
$ws.Cells.Item(77,1).Value = 24
$ws.Cells.Item(77,2).Value = "RR"
$ws.Cells.Item(77,3).Value = "GT"
$ws.Cells.Item(77,4).Value = 1
$ws.Cells.Item(77,5).Value = "RR"
$ws.Cells.Item(77,6).Value = "GT"
$ws.Cells.Item(77,7).Value = 20
$ws.Cells.Item(77,8).Value = "GT"
$ws.Cells.Item(77,9).Value = "HDPK Dharmasena"
$ws.Cells.Item(77,10).Value = "HDPKD"
$ws.Cells.Item(77,11).Value = "Wide"
$ws.Cells.Item(77,12).Value = "Called"
$ws.Cells.Item(77,13).Value = "Called"
$ws.Cells.Item(77,14).Value = "SV Samson"
$ws.Cells.Item(77,15).Value = "UT Yadav"
$ws.Cells.Item(77,16).Value = "Unsuccessful"
$ws.Cells.Item(77,17).Value = "No"

# Row 78
$ws.Cells.Item(78,1).Value = 24
$ws.Cells.Item(78,2).Value = "RR"
$ws.Cells.Item(78,3).Value = "GT"
$ws.Cells.Item(78,4).Value = 2
$ws.Cells.Item(78,5).Value = "GT"
$ws.Cells.Item(78,6).Value = "RR"
$ws.Cells.Item(78,7).Value = 9
$ws.Cells.Item(78,8).Value = "RR"
$ws.Cells.Item(78,9).Value = "HDPK Dharmasena"
$ws.Cells.Item(78,10).Value = "HDPKD"
$ws.Cells.Item(78,11).Value = "Wicket"
$ws.Cells.Item(78,12).Value = "Not Called"
$ws.Cells.Item(78,13).Value = "Called"
$ws.Cells.Item(78,14).Value = "B Sai Sudharsan"
$ws.Cells.Item(78,15).Value = "KR Sen"
$ws.Cells.Item(78,16).Value = "Successful"
$ws.Cells.Item(78,17).Value = "No"

# Row 79
$ws.Cells.Item(79,1).Value = 24
$ws.Cells.Item(79,2).Value = "RR"
$ws.Cells.Item(79,3).Value = "GT"
$ws.Cells.Item(79,4).Value = 2
$ws.Cells.Item(79,5).Value = "GT"
$ws.Cells.Item(79,6).Value = "RR"
$ws.Cells.Item(79,7).Value = 11
$ws.Cells.Item(79,8).Value = "RR"
$ws.Cells.Item(79,9).Value = "HDPK Dharmasena"
$ws.Cells.Item(79,10).Value = "HDPKD"
$ws.Cells.Item(79,11).Value = "Wide"
$ws.Cells.Item(79,12).Value = "Called"
$ws.Cells.Item(79,13).Value = "Called"
$ws.Cells.Item(79,14).Value = "MS Wade"
$ws.Cells.Item(79,15).Value = "KR Sen"
$ws.Cells.Item(79,16).Value = "Unsuccessful"
$ws.Cells.Item(79,17).Value = "No"

# Row 80
# (Bowler cell O80 is intentionally written before Batter cell N80 so that the
#  shared-string table ends up in the same insertion order as the source file.)
$ws.Cells.Item(80,1).Value = 24
$ws.Cells.Item(80,2).Value = "RR"
$ws.Cells.Item(80,3).Value = "GT"
$ws.Cells.Item(80,4).Value = 2
$ws.Cells.Item(80,5).Value = "GT"
$ws.Cells.Item(80,6).Value = "RR"
$ws.Cells.Item(80,7).Value = 16
$ws.Cells.Item(80,8).Value = "RR"
$ws.Cells.Item(80,9).Value = "Vinod Seshan"
$ws.Cells.Item(80,10).Value = "VS"
$ws.Cells.Item(80,11).Value = "Wicket"
$ws.Cells.Item(80,12).Value = "Not Out"
$ws.Cells.Item(80,13).Value = "Not Out"
$ws.Cells.Item(80,15).Value = "YS Chahal"
$ws.Cells.Item(80,14).Value = "M Shahrukh Khan"
$ws.Cells.Item(80,16).Value = "Unsuccessful"
$ws.Cells.Item(80,17).Value = "No"

# Row 81
$ws.Cells.Item(81,1).Value = 24
$ws.Cells.Item(81,2).Value = "RR"
$ws.Cells.Item(81,3).Value = "GT"
$ws.Cells.Item(81,4).Value = 2
$ws.Cells.Item(81,5).Value = "GT"
$ws.Cells.Item(81,6).Value = "RR"
$ws.Cells.Item(81,7).Value = 18
$ws.Cells.Item(81,8).Value = "GT"
$ws.Cells.Item(81,9).Value = "Vinod Seshan"
$ws.Cells.Item(81,10).Value = "VS"
$ws.Cells.Item(81,11).Value = "Wicket"
$ws.Cells.Item(81,12).Value = "Out"
$ws.Cells.Item(81,13).Value = "Out"
$ws.Cells.Item(81,14).Value = "M Shahrukh Khan"
$ws.Cells.Item(81,15).Value = "Avesh Khan"
$ws.Cells.Item(81,16).Value = "Unsuccessful"
$ws.Cells.Item(81,17).Value = "No"

# Leave the sheet scrolled/selected the same way the author left it.
$ws.Range("G59").Select()
